$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 13 through 30 (they are removed from the data range)
$ws.Range("A13:B30").EntireRow.Delete()

# Update the date values in column B for rows 2 through 12
# (values are entered with a leading apostrophe so Excel keeps them as
# literal text instead of auto-converting to date serials, then
# ClearFormats() drops the "quote prefix" style Excel applies so the
# cells keep the workbook's original/default formatting)
$ws.Range("B2").Value = "'2005-03-22"
$ws.Range("B3").Value = "'2005-03-30"
$ws.Range("B4").Value = "'2005-03-31"
$ws.Range("B5").Value = "'2005-11-02"
$ws.Range("B6").Value = "'2005-11-16"
$ws.Range("B7").Value = "'2006-01-09"
$ws.Range("B8").Value = "'2006-01-10"
$ws.Range("B9").Value = "'2006-01-11"
$ws.Range("B10").Value = "'2006-01-20"
$ws.Range("B11").Value = "'2006-01-23"
$ws.Range("B12").Value = "'2006-01-24"

$ws.Range("B2:B12").ClearFormats()
